# Actualización automática 2025-07-31 14:28:00
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M6").Value = 1041.33
$ws1.Range("I7").Value = 462.6
$ws1.Range("M7").Value = 1394.01
$ws1.Range("D15").Value = 950.4
$ws1.Range("L41").Value = 1696.17
$ws1.Range("D56").Value = "3 de 54"
$ws1.Range("M56").Value = "17 de 54"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F6").Value = 1041.33
$ws2.Range("F7").Value = 3246.21
$ws2.Range("F15").Value = 4931.71
$ws2.Range("F41").Value = 4722.13
$ws2.Range("F56").Value = 87786.59

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 2358.72
$ws3.Range("E3").Value = 25098.2876
$ws3.Range("F3").Value = 0.08590593827129216

$ws3.Range("D8").Value = 1236.46
$ws3.Range("E8").Value = -236.46
$ws3.Range("F8").Value = 1.23646

$ws3.Range("D15").Value = 8839.24
$ws3.Range("E15").Value = 4660.76
$ws3.Range("F15").Value = 0.6547585185185185

$ws3.Range("D16").Value = 63830.97
$ws3.Range("E16").Value = -12004.51
$ws3.Range("F16").Value = 1.231628978710875

$ws3.Range("D19").Value = 87786.59
$ws3.Range("E19").Value = 25919.86064517915
$ws3.Range("F19").Value = 0.7720458206363152
